$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-26 from 45170 to 45174
$ws.Range("C2:C26").Value = 45174
